$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value, even if it looks numeric
# (protects values like "318.23", "1.00", "19.00" from being coerced to
#  floating point numbers and losing their original text formatting).
function Set-TextValue {
    param($ws, $cellRef, $text)
    $cell = $ws.Range($cellRef)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.Formula = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# --- Price (column D) updates ---
Set-TextValue $ws "D2" "42.901.06"
Set-TextValue $ws "D3" "2.534.63"
Set-TextValue $ws "D5" "318.23"
Set-TextValue $ws "D6" "95.98"
Set-TextValue $ws "D7" "0.581"
Set-TextValue $ws "D10" "36.36"
Set-TextValue $ws "D12" "7.63"
Set-TextValue $ws "D14" "2.924.88"
Set-TextValue $ws "D15" "2.548.61"
Set-TextValue $ws "D16" "15.35"
Set-TextValue $ws "D17" "0.848"
Set-TextValue $ws "D18" "42.940.05"
Set-TextValue $ws "D19" "13.05"
Set-TextValue $ws "D21" "0.0₃0968"
Set-TextValue $ws "D22" "70.16"
Set-TextValue $ws "D23" "252.22"
Set-TextValue $ws "D24" "2.98"
Set-TextValue $ws "D26" "27.05"
Set-TextValue $ws "D29" "39.85"
Set-TextValue $ws "D31" "6.07"
Set-TextValue $ws "D32" "153.87"
Set-TextValue $ws "D33" "2.12"
Set-TextValue $ws "D34" "3.32"
Set-TextValue $ws "D35" "19.00"
Set-TextValue $ws "D40" "23.63"
Set-TextValue $ws "D41" "2.28"
Set-TextValue $ws "D46" "2.015.20"
Set-TextValue $ws "D47" "85.68"
Set-TextValue $ws "D49" "2.778.55"
Set-TextValue $ws "D50" "74.27"
Set-TextValue $ws "D51" "102.38"

# --- Volume(1h) (column E) updates ---
Set-TextValue $ws "E2" "  +0.98%  "
Set-TextValue $ws "E3" "  +0.88%  "
Set-TextValue $ws "E4" "  +0.00%  "
Set-TextValue $ws "E5" "  +4.72%  "
Set-TextValue $ws "E6" "  -0.25%  "
Set-TextValue $ws "E7" "  +0.58%  "
Set-TextValue $ws "E8" "  -0.05%  "
Set-TextValue $ws "E9" "  -0.72%  "
Set-TextValue $ws "E10" "  -0.54%  "
Set-TextValue $ws "E11" "  +0.35%  "
Set-TextValue $ws "E12" "  -0.33%  "
Set-TextValue $ws "E13" "  -0.29%  "
Set-TextValue $ws "E14" "  +0.73%  "
Set-TextValue $ws "E15" "  +0.12%  "
Set-TextValue $ws "E16" "  +2.51%  "
Set-TextValue $ws "E17" "  -1.10%  "
Set-TextValue $ws "E18" "  +0.62%  "
Set-TextValue $ws "E19" "  +1.14%  "
Set-TextValue $ws "E20" "  +3.43%  "
Set-TextValue $ws "E21" "  -0.26%  "
Set-TextValue $ws "E22" "  -1.14%  "
Set-TextValue $ws "E23" "  +0.89%  "
Set-TextValue $ws "E24" "  +2.28%  "
Set-TextValue $ws "E25" "  +0.31%  "
Set-TextValue $ws "E26" "  +1.06%  "
Set-TextValue $ws "E27" "  -0.17%  "
Set-TextValue $ws "E28" "  +4.18%  "
Set-TextValue $ws "E29" "  +5.32%  "
Set-TextValue $ws "E30" "  -0.25%  "
Set-TextValue $ws "E31" "  +2.15%  "
Set-TextValue $ws "E32" "  -1.66%  "
Set-TextValue $ws "E33" "  +3.39%  "
Set-TextValue $ws "E34" "  +2.46%  "
Set-TextValue $ws "E35" "  +3.60%  "
Set-TextValue $ws "E36" "  +0.48%  "
Set-TextValue $ws "E37" "  +0.08%  "
Set-TextValue $ws "E38" "  -3.36%  "
Set-TextValue $ws "E39" "  +0.10%  "
Set-TextValue $ws "E40" "  -5.86%  "
Set-TextValue $ws "E41" "  +9.37%  "
Set-TextValue $ws "E46" "  -0.57%  "
Set-TextValue $ws "E47" "  +1.00%  "
Set-TextValue $ws "E48" "  -1.96%  "
Set-TextValue $ws "E49" "  +0.51%  "
Set-TextValue $ws "E50" "  +3.31%  "
Set-TextValue $ws "E51" "  +1.09%  "

# --- Row reordering: rows 42-43 swap (RenderToken <-> VeChain),
#     rows 44-45 swap (NEARProtocol <-> FirstDigitalUSD), each with updated % values ---
Set-TextValue $ws "B42" "VeChain"
Set-TextValue $ws "C42" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D42" "0.0305"
Set-TextValue $ws "E42" "  +1.88%  "

Set-TextValue $ws "B43" "RenderToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D43" "3.81"
Set-TextValue $ws "E43" "  -0.31%  "

Set-TextValue $ws "B44" "FirstDigitalUSD"
Set-TextValue $ws "C44" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D44" "1.00"
Set-TextValue $ws "E44" "  +0.56%  "

Set-TextValue $ws "B45" "NEARProtocol"
Set-TextValue $ws "C45" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D45" "3.33"
Set-TextValue $ws "E45" "  -1.25%  "
